# Daily IST report: add CSV/MD/XLSX
# Insert a new date column (2026-02-19) at column H, shifting
# total_files -> I and unique_days -> J, and recompute both summary
# columns for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "total_files" column (H).
$ws.Columns.Item(8).Insert()

# Match the width used by the other date columns (D:G are 12 wide).
$ws.Columns.Item(8).ColumnWidth = 11.17

# New header cell for the inserted date column. The leading apostrophe
# forces Excel to store it as literal text instead of auto-converting
# the date-shaped string into a date value (matching the other date
# headers in row 1, which are plain text too). Re-apply G1's format
# afterwards so H1 keeps the same plain bold header style as the rest
# of the date columns instead of the "quote prefix" variant.
$ws.Cells.Item(1, 8).Value = "'2026-02-19"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# total_files / unique_days headers now live one column to the right.
$ws.Cells.Item(1, 9).Value = "total_files"
$ws.Cells.Item(1, 10).Value = "unique_days"

# Per-row data: new 2026-02-19 submission count (H), recomputed
# total_files (I) and recomputed unique_days (J).
$rowData = @(
    "1,2,2",
    "0,1,1",
    "1,1,1",
    "1,2,2",
    "0,0,0",
    "1,2,2",
    "1,2,2",
    "1,2,2",
    "1,2,2",
    "1,1,1",
    "1,2,2",
    "1,2,2",
    "1,2,2",
    "0,1,1",
    "0,0,0",
    "1,2,2",
    "1,2,2",
    "1,2,2",
    "0,1,1",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "0,1,1",
    "1,2,2",
    "0,0,0",
    "0,3,1",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "1,2,2",
    "1,2,2",
    "1,2,2",
    "1,2,2",
    "1,2,2",
    "1,1,1",
    "0,1,1",
    "0,0,0",
    "0,0,0",
    "1,2,2",
    "1,2,2",
    "1,2,2",
    "0,21,1",
    "0,0,0",
    "0,1,1",
    "1,2,2",
    "1,2,2",
    "0,0,0",
    "1,2,2",
    "1,1,1",
    "0,0,0",
    "1,2,2",
    "1,1,1",
    "0,0,0",
    "1,2,2",
    "0,0,0",
    "1,1,1",
    "0,1,1",
    "1,1,1",
    "1,2,2",
    "0,1,1",
    "0,0,0",
    "0,0,0",
    "1,2,2",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "1,2,2",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "1,14,2",
    "0,0,0",
    "0,1,1",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "1,1,1",
    "1,1,1",
    "1,2,2",
    "0,0,0",
    "0,1,1",
    "0,0,0",
    "0,1,1",
    "0,0,0",
    "1,1,1",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "1,1,1",
    "1,31,2",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "0,0,0",
    "0,1,1",
    "1,2,2",
    "0,0,0",
    "1,2,2",
    "0,12,1",
    "1,2,2",
    "1,2,2",
    "0,0,0",
    "0,1,1",
    "0,0,0"
)

$r = 2
foreach ($line in $rowData) {
    $parts = $line.Split(",")
    $ws.Cells.Item($r, 8).Value = [int]$parts[0]
    $ws.Cells.Item($r, 9).Value = [int]$parts[1]
    $ws.Cells.Item($r, 10).Value = [int]$parts[2]
    $r++
}
